# Relabel and relocate the single file-listing row: update the file name
# and its corresponding last-write-time to reflect the new artifact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Screenshot 2023-03-15 084217.png"
$ws.Range("C3").Value = "3/15/2023 8:42:17 AM"
